$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds "K" values; update per regenerated save_data
$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 4
$ws.Range("G4").Value = 3
$ws.Range("G5").Value = 8
